# Add season record columns (Wins, Losses, Ties) to the BOS_2021 roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers in row 1, matching the style used by the other header cells (AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Season record values for every data row (2-58): 92 wins, 70 losses, 0 ties.
$lastRow = 58
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 92  # AD
    $ws.Cells.Item($r, 31).Value = 70  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
